$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GSM")
$ws.Rows.Item(1).Delete()
Write-Host "deleted"
